$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 4.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("AA2").Value = 21
$ws.Range("AW2").Value = 6
$ws.Range("BA2").Value = 151

# Row 3
$ws.Range("H3").Value = 3.8
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("AE3").Value = 21
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 51
$ws.Range("AK3").Value = 41
$ws.Range("AZ3").Value = 126

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("W4").Value = 4.75
$ws.Range("X4").Value = 6.5
$ws.Range("AG4").Value = 10
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 67
$ws.Range("AL4").Value = 67
$ws.Range("AX4").Value = 34

# Row 6
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 4.33
$ws.Range("I6").Value = 1.42
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 2.6
$ws.Range("L6").Value = 1.91
$ws.Range("U6").Value = 1.62
$ws.Range("V6").Value = 2.2
$ws.Range("W6").Value = 23
$ws.Range("X6").Value = 41
$ws.Range("Y6").Value = 21
$ws.Range("Z6").Value = 67
$ws.Range("AA6").Value = 41
$ws.Range("AB6").Value = 41
$ws.Range("AD6").Value = 9.5
$ws.Range("AE6").Value = 15
$ws.Range("AH6").Value = 9
$ws.Range("AI6").Value = 8.5
$ws.Range("AJ6").Value = 11
$ws.Range("AK6").Value = 11
$ws.Range("AN6").Value = 8
$ws.Range("AO6").Value = 26
$ws.Range("AP6").Value = 29
$ws.Range("AQ6").Value = 81
$ws.Range("AS6").Value = 151
$ws.Range("AU6").Value = 8
$ws.Range("AW6").Value = 3.75
$ws.Range("AX6").Value = 7
$ws.Range("AZ6").Value = 17

# Row 7
$ws.Range("G7").Value = 2.45
$ws.Range("I7").Value = 2.7
$ws.Range("J7").Value = 3
$ws.Range("L7").Value = 3.2
$ws.Range("Q7").Value = 1.83
$ws.Range("R7").Value = 1.98
$ws.Range("X7").Value = 13
$ws.Range("Y7").Value = 10
$ws.Range("Z7").Value = 23
$ws.Range("AG7").Value = 10
$ws.Range("AJ7").Value = 26
$ws.Range("AS7").Value = 151
$ws.Range("AW7").Value = 4.75
$ws.Range("AZ7").Value = 41

# Row 8
$ws.Range("G8").Value = 2.7
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 2.45
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.33
$ws.Range("T8").Value = 3.25
$ws.Range("W8").Value = 11
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 10
$ws.Range("AC8").Value = 13
$ws.Range("AH8").Value = 13
$ws.Range("AI8").Value = 9.5
$ws.Range("AJ8").Value = 23
$ws.Range("AK8").Value = 19
$ws.Range("AL8").Value = 23
$ws.Range("AM8").Value = 151
$ws.Range("AQ8").Value = 41
$ws.Range("AT8").Value = 3.25
$ws.Range("AW8").Value = 4.75
$ws.Range("AX8").Value = 13
$ws.Range("AY8").Value = 21
$ws.Range("BB8").Value = 126
$ws.Range("BC8").Value = 451
